# Updated cryptos list - apply new Price (D) and Volume(1h) (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.744.61"
$ws.Range("D2").Style = $ws.Range("C2").Style
$ws.Range("E2").Value = "  +0.29%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.291.30"
$ws.Range("D3").Style = $ws.Range("C3").Style
$ws.Range("E3").Value = "  +0.02%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "113.64"
$ws.Range("D5").Style = $ws.Range("C5").Style
$ws.Range("E5").Value = "  +17.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "268.12"
$ws.Range("D6").Style = $ws.Range("C6").Style
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  +0.22%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.620"
$ws.Range("D9").Style = $ws.Range("C9").Style
$ws.Range("E9").Value = "  +1.49%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.00"
$ws.Range("D10").Style = $ws.Range("C10").Style
$ws.Range("E10").Value = "  +4.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0937"
$ws.Range("D11").Style = $ws.Range("C11").Style

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.78"
$ws.Range("D12").Style = $ws.Range("C12").Style
$ws.Range("E12").Value = "  +12.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.107"
$ws.Range("D13").Style = $ws.Range("C13").Style
$ws.Range("E13").Value = "  +0.99%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.60"
$ws.Range("D14").Style = $ws.Range("C14").Style
$ws.Range("E14").Value = "  +2.96%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.635.10"
$ws.Range("D15").Style = $ws.Range("C15").Style
$ws.Range("E15").Value = "  +0.02%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.844"
$ws.Range("D16").Style = $ws.Range("C16").Style
$ws.Range("E16").Value = "  -0.46%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.291.47"
$ws.Range("D17").Style = $ws.Range("C17").Style
$ws.Range("E17").Value = "  +0.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.629.61"
$ws.Range("D18").Style = $ws.Range("C18").Style
$ws.Range("E18").Value = "  +0.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000109"
$ws.Range("D19").Style = $ws.Range("C19").Style
$ws.Range("E19").Value = "  +1.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.52"
$ws.Range("D20").Style = $ws.Range("C20").Style
$ws.Range("E20").Value = "  +4.73%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.55"
$ws.Range("D21").Style = $ws.Range("C21").Style
$ws.Range("E21").Value = "  +0.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.46"
$ws.Range("D22").Style = $ws.Range("C22").Style
$ws.Range("E22").Value = "  -0.84%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.75"
$ws.Range("D23").Style = $ws.Range("C23").Style
$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.78"
$ws.Range("D24").Style = $ws.Range("C24").Style
$ws.Range("E24").Value = "  +6.29%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.81"
$ws.Range("D25").Style = $ws.Range("C25").Style
$ws.Range("E25").Value = "  +10.46%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = $ws.Range("C26").Style
$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.61"
$ws.Range("D27").Style = $ws.Range("C27").Style
$ws.Range("E27").Value = "  +3.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "41.92"
$ws.Range("D28").Style = $ws.Range("C28").Style
$ws.Range("E28").Value = "  +3.80%  "

$ws.Range("E29").Value = "  -2.19%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.26"
$ws.Range("D30").Style = $ws.Range("C30").Style
$ws.Range("E30").Value = "  -0.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "176.52"
$ws.Range("D31").Style = $ws.Range("C31").Style
$ws.Range("E31").Value = "  +0.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0928"
$ws.Range("D32").Style = $ws.Range("C32").Style
$ws.Range("E32").Value = "  +4.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.56"
$ws.Range("D33").Style = $ws.Range("C33").Style
$ws.Range("E33").Value = "  -1.27%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.55"
$ws.Range("D34").Style = $ws.Range("C34").Style
$ws.Range("E34").Value = "  +3.31%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.126"
$ws.Range("D35").Style = $ws.Range("C35").Style
$ws.Range("E35").Value = "  +0.12%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.74"
$ws.Range("D36").Style = $ws.Range("C36").Style
$ws.Range("E36").Value = "  +8.79%  "

$ws.Range("E37").Value = "  +0.50%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0357"
$ws.Range("D38").Style = $ws.Range("C38").Style
$ws.Range("E38").Value = "  +0.74%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.82"
$ws.Range("D39").Style = $ws.Range("C39").Style
$ws.Range("E39").Value = "  +12.50%  "

$ws.Range("E40").Value = "  +5.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.88"
$ws.Range("D41").Style = $ws.Range("C41").Style
$ws.Range("E41").Value = "  +12.68%  "

$ws.Range("E42").Value = "  +2.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.41"
$ws.Range("D43").Style = $ws.Range("C43").Style
$ws.Range("E43").Value = "  +11.97%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.13"
$ws.Range("D44").Style = $ws.Range("C44").Style
$ws.Range("E44").Value = "  +17.50%  "

$ws.Range("E45").Value = "  +6.03%  "

$ws.Range("E46").Value = "  +0.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.72"
$ws.Range("D47").Style = $ws.Range("C47").Style
$ws.Range("E47").Value = "  -0.68%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.25"
$ws.Range("D48").Style = $ws.Range("C48").Style
$ws.Range("E48").Value = "  +5.84%  "

$ws.Range("E49").Value = "  -1.43%  "

$ws.Range("E50").Value = "  +3.37%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.449"
$ws.Range("D51").Style = $ws.Range("C51").Style
$ws.Range("E51").Value = "  +5.03%  "
